$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.02739999999999
$ws.Range("A10").Value = -21.47699999999998
$ws.Range("A12").Value = -21.59249999999999
$ws.Range("C12").Value = -12.6498
$ws.Range("D12").Value = -8.540800000000004
$ws.Range("D13").Value = -9.009299999999984
$ws.Range("C17").Value = -13.48969999999999
$ws.Range("A18").Value = -22.1127
$ws.Range("D21").Value = -7.898699999999995
$ws.Range("C26").Value = -12.68130000000001
$ws.Range("C27").Value = -12.89789999999999
$ws.Range("C28").Value = -13.91639999999999
$ws.Range("D36").Value = -7.382799999999996
$ws.Range("A37").Value = -20.10079999999998
$ws.Range("C37").Value = -13.17600000000001
$ws.Range("D38").Value = -8.212699999999998
$ws.Range("D41").Value = -8.055899999999996
$ws.Range("D52").Value = -7.848500000000004
$ws.Range("A55").Value = -22.2159
$ws.Range("C65").Value = -12.25549999999999
$ws.Range("D67").Value = -7.347899999999995
$ws.Range("A68").Value = -21.44789999999999
$ws.Range("C73").Value = -11.33550000000001
$ws.Range("A77").Value = -20.16749999999998
$ws.Range("A78").Value = -20.22459999999999
$ws.Range("C84").Value = -12.96219999999999
$ws.Range("C85").Value = -13.7773
$ws.Range("D89").Value = -8.267499999999997
$ws.Range("C93").Value = -10.2499
$ws.Range("C95").Value = -12.9498
$ws.Range("D95").Value = -7.484900000000006
$ws.Range("C98").Value = -12.53510000000001
$ws.Range("C99").Value = -12.4516
$ws.Range("C101").Value = -13.52580000000001
$ws.Range("D105").Value = -8.349900000000003
